$wb = $excel.ActiveWorkbook

# F-column ("想去人数") updates that apply identically to both the
# "展览" and "全部类型" worksheets.
$updates = @{
    2  = 1474
    5  = 2201
    7  = 1352
    9  = 140
    11 = 323
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
